$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7025729417800903
$ws.Range("B1").Value = 1.383336305618286
$ws.Range("C1").Value = 5.068563461303711
$ws.Range("D1").Value = 1.604281067848206
$ws.Range("E1").Value = 0.9150964617729187
